$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 390088.6
$ws.Range("I33").Value = 399.85715
$ws.Range("J33").Value = 1299362.4
$ws.Range("K33").Value = 399.85715
$ws.Range("L33").Value = 1299362.4
$ws.Range("M33").Value = -170.85715
$ws.Range("N33").Value = -1299820.4
# Row 64
$ws.Range("H64").Value = 41313.19
$ws.Range("I64").Value = 113665.336
$ws.Range("J64").Value = 3009.1177
$ws.Range("K64").Value = 113665.336
$ws.Range("L64").Value = 3009.1177
$ws.Range("M64").Value = -113417.336
$ws.Range("N64").Value = -3505.1177
# Row 67
$ws.Range("H67").Value = 41313.19
$ws.Range("I67").Value = 113665.336
$ws.Range("J67").Value = 3009.1177
$ws.Range("K67").Value = 113665.336
$ws.Range("L67").Value = 3009.1177
$ws.Range("M67").Value = -112807.336
$ws.Range("N67").Value = -4725.1177
# Row 74
$ws.Range("H74").Value = 3899.8572
$ws.Range("I74").Value = 3324.75
$ws.Range("K74").Value = 3324.75
$ws.Range("M74").Value = -2388.75
# Row 76
$ws.Range("H76").Value = 4749.25
$ws.Range("I76").Value = 4569.857
$ws.Range("J76").Value = 5000.4
$ws.Range("K76").Value = 4569.857
$ws.Range("L76").Value = 5000.4
$ws.Range("M76").Value = -4254.857
$ws.Range("N76").Value = -5630.4
# Row 77
$ws.Range("H77").Value = 3899.8572
$ws.Range("I77").Value = 3324.75
$ws.Range("K77").Value = 16623.75
$ws.Range("M77").Value = -11943.75
# Row 79
$ws.Range("H79").Value = 4749.25
$ws.Range("I79").Value = 4569.857
$ws.Range("J79").Value = 5000.4
$ws.Range("K79").Value = 4569.857
$ws.Range("L79").Value = 5000.4
$ws.Range("M79").Value = -3477.857
$ws.Range("N79").Value = -7184.4
# Row 113
$ws.Range("H113").Value = 57453.332
$ws.Range("I113").Value = 126921.75
$ws.Range("J113").Value = 1878.6
$ws.Range("K113").Value = 126921.75
$ws.Range("L113").Value = 1878.6
$ws.Range("M113").Value = -123667.75
$ws.Range("N113").Value = -8386.6
# Row 123
$ws.Range("H123").Value = 28893
$ws.Range("J123").Value = 28893
$ws.Range("L123").Value = 28893
$ws.Range("N123").Value = -38693
# Row 135
$ws.Range("H135").Value = 1685.5405
$ws.Range("I135").Value = 727.4167
$ws.Range("J135").Value = 3454.3845
$ws.Range("K135").Value = 6546.7503
$ws.Range("L135").Value = 31089.4605
$ws.Range("M135").Value = -4011.7503
$ws.Range("N135").Value = -36159.4605

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 31765.924
$ws.Range("I32").Value = 9363.5
$ws.Range("K32").Value = 9363.5
$ws.Range("M32").Value = -9076.5
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
# Row 61
$ws.Range("H61").Value = 1901.3846
$ws.Range("I61").Value = 1599.091
$ws.Range("J61").Value = 3564
$ws.Range("K61").Value = 1599.091
$ws.Range("L61").Value = 3564
$ws.Range("M61").Value = -1387.091
$ws.Range("N61").Value = -3988
# Row 74
$ws.Range("H74").Value = 918.8095
$ws.Range("I74").Value = 898.7059
$ws.Range("J74").Value = 1004.25
$ws.Range("K74").Value = 898.7059
$ws.Range("L74").Value = 1004.25
$ws.Range("M74").Value = -24.70590000000004
$ws.Range("N74").Value = -2752.25
# Row 77
$ws.Range("H77").Value = 918.8095
$ws.Range("I77").Value = 898.7059
$ws.Range("J77").Value = 1004.25
$ws.Range("K77").Value = 4493.529500000001
$ws.Range("L77").Value = 5021.25
$ws.Range("M77").Value = -125.5295000000006
$ws.Range("N77").Value = -13757.25
# Row 88
$ws.Range("H88").Value = 1000
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -594
$ws.Range("N88").ClearContents()
# Row 91
$ws.Range("H91").Value = 1000
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 404
$ws.Range("N91").ClearContents()
# Row 122
$ws.Range("H122").Value = 2573.2
$ws.Range("I122").Value = 2100.2666
$ws.Range("J122").Value = 3992
$ws.Range("K122").Value = 6300.7998
$ws.Range("L122").Value = 11976
$ws.Range("M122").Value = -3850.7998
$ws.Range("N122").Value = -16876
# Row 136
$ws.Range("H136").Value = 1901.3846
$ws.Range("I136").Value = 1599.091
$ws.Range("J136").Value = 3564
$ws.Range("K136").Value = 4797.272999999999
$ws.Range("L136").Value = 10692
$ws.Range("M136").Value = -2247.272999999999
$ws.Range("N136").Value = -15792

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 9
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
# Row 44
$ws.Range("H44").Value = 11225
$ws.Range("J44").Value = 13266.667
$ws.Range("L44").Value = 13266.667
$ws.Range("N44").Value = -14260.667
# Row 107
$ws.Range("H107").Value = 55581800
$ws.Range("I107").Value = 58851228
$ws.Range("K107").Value = 58851228
$ws.Range("M107").Value = -58849308

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1766.1724
$ws.Range("I58").Value = 1602.6666
$ws.Range("J58").Value = 2033.7273
$ws.Range("K58").Value = 1602.6666
$ws.Range("L58").Value = 2033.7273
$ws.Range("M58").Value = -1399.6666
$ws.Range("N58").Value = -2439.7273
# Row 62
$ws.Range("H62").Value = 2657.1428
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
# Row 65
$ws.Range("H65").Value = 2657.1428
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380
# Row 136
$ws.Range("H136").Value = 1766.1724
$ws.Range("I136").Value = 1602.6666
$ws.Range("J136").Value = 2033.7273
$ws.Range("K136").Value = 4807.9998
$ws.Range("L136").Value = 6101.1819
$ws.Range("M136").Value = -2257.9998
$ws.Range("N136").Value = -11201.1819

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 102
$ws.Range("H102").Value = 4832.6665
$ws.Range("J102").Value = 4999
$ws.Range("L102").Value = 14997
$ws.Range("N102").Value = -19865
# Row 113
$ws.Range("H113").Value = 925.9655
$ws.Range("J113").Value = 613.9474
$ws.Range("L113").Value = 1841.8422
$ws.Range("N113").Value = -6181.8422
# Row 131
$ws.Range("H131").Value = 846993.4399999999
$ws.Range("I131").Value = 602
$ws.Range("J131").Value = 1012952.56
$ws.Range("K131").Value = 1806
$ws.Range("L131").Value = 3038857.68
$ws.Range("M131").Value = 3234
$ws.Range("N131").Value = -3048937.68

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 62940.688
$ws.Range("I70").Value = 96240.82000000001
$ws.Range("J70").Value = 6586.615
$ws.Range("K70").Value = 96240.82000000001
$ws.Range("L70").Value = 6586.615
$ws.Range("M70").Value = -95970.82000000001
$ws.Range("N70").Value = -7126.615
# Row 73
$ws.Range("H73").Value = 62940.688
$ws.Range("I73").Value = 96240.82000000001
$ws.Range("J73").Value = 6586.615
$ws.Range("K73").Value = 96240.82000000001
$ws.Range("L73").Value = 6586.615
$ws.Range("M73").Value = -95304.82000000001
$ws.Range("N73").Value = -8458.615
# Row 102
$ws.Range("H102").Value = 2452.6316
$ws.Range("I102").Value = 2625.375
$ws.Range("J102").Value = 2327
$ws.Range("K102").Value = 2625.375
$ws.Range("L102").Value = 2327
$ws.Range("M102").Value = -1003.375
$ws.Range("N102").Value = -5571

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 2645.0605
$ws.Range("I122").Value = 2605.25
$ws.Range("J122").Value = 2706.3076
$ws.Range("K122").Value = 7815.75
$ws.Range("L122").Value = 8118.9228
$ws.Range("M122").Value = -5365.75
$ws.Range("N122").Value = -13018.9228

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2934.7334
$ws.Range("I122").Value = 2009.6666
$ws.Range("J122").Value = 6635
$ws.Range("K122").Value = 6028.9998
$ws.Range("L122").Value = 19905
$ws.Range("M122").Value = -3578.9998
$ws.Range("N122").Value = -24805
# Row 126
$ws.Range("H126").Value = 1450.25
$ws.Range("I126").Value = 1394.7222
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 4184.1666
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -1714.1666
$ws.Range("N126").Value = -10790
# Row 132
$ws.Range("H132").Value = 9896.5
$ws.Range("I132").Value = 5922.846
$ws.Range("J132").Value = 17276.143
$ws.Range("K132").Value = 17768.538
$ws.Range("L132").Value = 51828.429
$ws.Range("M132").Value = -15238.538
$ws.Range("N132").Value = -56888.429
# Row 136
$ws.Range("H136").Value = 29695.475
$ws.Range("I136").Value = 125888
$ws.Range("J136").Value = 5647.3438
$ws.Range("K136").Value = 377664
$ws.Range("L136").Value = 16942.0314
$ws.Range("M136").Value = -375114
$ws.Range("N136").Value = -22042.0314
